$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.426.79'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.619.33'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.94'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0609'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.21'
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("D12").Value = '1.846.79'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '1.620.44'
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.81'
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '236.60'
$ws.Range("E17").Value = '  +8.18%  '
$ws.Range("D18").Value = '26.428.45'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.84'
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("D20").Value = '0.0₃0726'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.19'
$ws.Range("E23").Value = '  +2.99%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.12'
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("E25").Value = '  +1.55%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.07'
$ws.Range("E27").Value = '  +0.73%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("E29").Value = '  +2.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("E31").Value = '  -0.57%  '
$ws.Range("D32").Value = '1.517.16'
$ws.Range("E32").Value = '  +5.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.26'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.98'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.50'
$ws.Range("E35").Value = '  +1.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.569'
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("E42").Value = '  +0.77%  '
$ws.Range("D43").Value = '1.759.70'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.763'
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.04'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.911'
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.33'
$ws.Range("E47").Value = '  +2.76%  '
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0963'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("E51").Value = '  +0.86%  '

Write-Output "done"
